$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add value 5 into newly-used cells in columns X, Y, Z for several rows.
$ws.Range("X3").Value = 5
$ws.Range("Y3").Value = 5
$ws.Range("Z3").Value = 5

$ws.Range("X4").Value = 5
$ws.Range("Y4").Value = 5
$ws.Range("Z4").Value = 5

$ws.Range("X6").Value = 5

$ws.Range("X8").Value = 5
$ws.Range("Y8").Value = 5
$ws.Range("Z8").Value = 5

$ws.Range("X9").Value = 5
$ws.Range("Y9").Value = 5
$ws.Range("Z9").Value = 5

$ws.Range("X11").Value = 5
$ws.Range("Z11").Value = 5

$ws.Range("Y13").Value = 5
$ws.Range("Z13").Value = 5

$ws.Range("X20").Value = 5
$ws.Range("Y20").Value = 5
$ws.Range("Z20").Value = 5

$ws.Range("X25").Value = 5
$ws.Range("Y25").Value = 5
$ws.Range("Z25").Value = 5

$ws.Range("Y26").Value = 5
$ws.Range("Z26").Value = 5

# Move the active selection to X3, matching the saved sheet view state.
$ws.Range("X3").Select()
